$wb = $excel.ActiveWorkbook

# --- Teams sheet: change the date for the last match day (rows 70-73) ---
# 30/08/2025 (serial 45899) -> 31/08/2025 (serial 45900)
$teams = $wb.Worksheets.Item("Teams")
$teams.Range("A70").Value = 45900
$teams.Range("A71").Value = 45900
$teams.Range("A72").Value = 45900
$teams.Range("A73").Value = 45900

# Update the active selection on the Teams sheet to A73
$teams.Range("A73").Select()

# --- Stats sheet: scroll the view down a bit (topLeftCell A455 -> A458) ---
$stats = $wb.Worksheets.Item("Stats")
$stats.Activate()
$stats.Application.ActiveWindow.ScrollRow = 458
$stats.Range("D466").Select()
